$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 47, shifting the existing row 47 (and below) down to row 48.
$ws.Rows.Item(47).Insert()

# Apply the same date number format used by the rest of column D to the new D47 cell.
$ws.Cells.Item(47, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 47 with the weekly update values.
$ws.Cells.Item(47, 1).Value = 3
$ws.Cells.Item(47, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = 44746
$ws.Cells.Item(47, 5).Value = 5
$ws.Cells.Item(47, 6).Value = 100112035
$ws.Cells.Item(47, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 103
$ws.Cells.Item(47, 11).Value = 15000
$ws.Cells.Item(47, 12).Value = 16000
$ws.Cells.Item(47, 13).Value = 15563
$ws.Cells.Item(47, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(47, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(47, 16).Value = 1038
$ws.Cells.Item(47, 17).Value = 15
$ws.Cells.Item(47, 18).Value = "Hortaliza"
